$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), reusing the formatting from the
# existing header cell H1 (bold font, border, centered alignment) so the
# new cells share the same style index as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I and J, rows 2-11
$values = @{
    2  = @(1, 2)
    3  = @(1, 4)
    4  = @(1, 6)
    5  = @(6, 6)
    6  = @(1, 4)
    7  = @(5, 7)
    8  = @(8, 8)
    9  = @(4, 4)
    10 = @(1, 1)
    11 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
